$d = $word.ActiveDocument

$replacements = @(
    @{old="36×49=1764"; new="29×81=2349"},
    @{old="21×30=630"; new="90×70=6300"},
    @{old="27×60=1620"; new="69×48=3312"},
    @{old="73×99=7227"; new="83×87=7221"},
    @{old="67×77=5159"; new="23×96=2208"},
    @{old="50×14=700"; new="32×94=3008"},
    @{old="20×99=1980"; new="79×79=6241"},
    @{old="28×79=2212"; new="92×20=1840"},
    @{old="62×85=5270"; new="96×67=6432"},
    @{old="86×48=4128"; new="50×85=4250"},
    @{old="49×81=3969"; new="46×63=2898"},
    @{old="45×61=2745"; new="90×18=1620"},
    @{old="29×18=522"; new="67×23=1541"},
    @{old="58×95=5510"; new="35×54=1890"},
    @{old="69×22=1518"; new="95×20=1900"},
    @{old="55×36=1980"; new="11×34=374"},
    @{old="53×81=4293"; new="95×45=4275"},
    @{old="68×36=2448"; new="72×23=1656"},
    @{old="54×42=2268"; new="76×77=5852"},
    @{old="78×70=5460"; new="33×43=1419"},
    @{old="49×41=2009"; new="62×47=2914"},
    @{old="12×65=780"; new="89×81=7209"},
    @{old="14×36=504"; new="69×52=3588"},
    @{old="64×65=4160"; new="46×96=4416"},
    @{old="37×63=2331"; new="42×11=462"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
